# Week 22 group exercise: add a new sub-bullet under the "A /<year> POST
# endpoint ..." item, describing that a form (or the API url) can be used
# to pass data in. The new bullet sits at the same outline level as the
# existing "Make sure you see your updates when you go to /all" item.

$d = $word.ActiveDocument

# Locate the paragraph whose text is the POST-endpoint bullet; the new
# paragraph must be inserted immediately after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*POST endpoint that lets you add additional data*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the POST-endpoint paragraph"
}

# Create a new (initially empty) paragraph right after the target one.
$target.Range.InsertParagraphAfter() | Out-Null
$newPara = $target.Next()

# The new bullet belongs one level deeper than the POST-endpoint bullet
# (ListLevelNumber is 1-based: level 3 == w:ilvl val="2"), matching the
# sibling "Make sure you see your updates..." bullet below it.
$newPara.Range.ListFormat.ListLevelNumber = 3

# Build the paragraph content as two separate runs via a WordOpenXML
# fragment so the text isn't merged into a single run.
$xml = @'
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>You can use a form or pass the data in to the API</w:t></w:r><w:r><w:t xml:space="preserve"> url</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$newPara.Range.InsertXML($xml) | Out-Null

Write-Host "Inserted new sub-bullet after the POST-endpoint item."
